$d = $word.ActiveDocument

# Locate the substring "Cục CSQLHC về TTXH" inside the run whose full text is
# "do Cục CSQLHC về TTXH cấp". Toggling a character formatting property on
# just that sub-range forces Word to split the parent run into three runs
# (before / middle / after) while leaving the run properties themselves
# unchanged, matching the target diff.
$rng = $d.Content
$found = $rng.Find.Execute("Cục CSQLHC về TTXH", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found) {
    $rng.Bold = 1
    $rng.Bold = 0
}
